# Add an IS_NUMBER column (H) to the column_names sheet that flags whether
# a row's DATA_TYPE (column D) is a numeric SQL type. This supports the new
# "functions ... for the sql fragment generator" that need to branch on
# whether a column is numeric.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Cells.Item(1, 8).Value = "IS_NUMBER"

$numericTypes = @("int", "smallint", "tinyint", "mediumint", "float", "decimal")

$lastRow = 239
for ($r = 2; $r -le $lastRow; $r++) {
    $dataType = $ws.Cells.Item($r, 4).Value2
    $isNumber = $numericTypes -contains $dataType
    $ws.Cells.Item($r, 8).Value = $isNumber
}

# Reset the view: scroll back to the top and move the active selection.
$ws.Range("I7").Select()
